# Apply the latest automation run's updates to the FAP report sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure every cell we touch stays plain text (the source data keeps
# leading zeros in CNPJ numbers and formatted strings like "1,0000"/"2025"),
# so force Text format across the existing table before writing any values.
$ws.Range("A1:I14").NumberFormat = "@"

# The "Estab_Nome" column (D) is no longer exported - remove it, which
# shifts UF/Municipio/Vigencia/Aliquota/Data_Consulta one column to the left.
$ws.Columns("D").Delete()

# Pre-format the two new rows as text as well before writing their values.
$ws.Range("A15:H16").NumberFormat = "@"

# --- Refresh CNPJ formatting (remove punctuation) and the query timestamps
# for the rows that were already present. Columns after the delete are:
# A=CNPJ_Raiz B=Razao_Social C=CNPJ_Estab D=UF E=Municipio F=Vigencia G=Aliquota H=Data_Consulta

$ws.Cells.Item(2,1).Value  = "53458313"
$ws.Cells.Item(2,3).Value  = "53458313000154"
$ws.Cells.Item(2,8).Value  = "31/10/2025 10:36:13"

$ws.Cells.Item(3,1).Value  = "02618143"
$ws.Cells.Item(3,3).Value  = "02618143000197"
$ws.Cells.Item(3,8).Value  = "31/10/2025 10:36:36"

$ws.Cells.Item(4,1).Value  = "43228913"
$ws.Cells.Item(4,3).Value  = "43228913000172"
$ws.Cells.Item(4,8).Value  = "31/10/2025 10:36:59"

$ws.Cells.Item(5,1).Value  = "03720943"
$ws.Cells.Item(5,3).Value  = "03720943000187"
$ws.Cells.Item(5,8).Value  = "31/10/2025 10:37:23"

$ws.Cells.Item(6,1).Value  = "03720943"
$ws.Cells.Item(6,3).Value  = "03720943000349"
$ws.Cells.Item(6,8).Value  = "31/10/2025 10:37:36"

$ws.Cells.Item(7,1).Value  = "03720943"
$ws.Cells.Item(7,3).Value  = "03720943000772"
$ws.Cells.Item(7,8).Value  = "31/10/2025 10:37:48"

$ws.Cells.Item(8,1).Value  = "03720943"
$ws.Cells.Item(8,3).Value  = "03720943001159"
$ws.Cells.Item(8,8).Value  = "31/10/2025 10:38:00"

$ws.Cells.Item(9,1).Value  = "03720943"
$ws.Cells.Item(9,3).Value  = "03720943001230"
$ws.Cells.Item(9,8).Value  = "31/10/2025 10:38:13"

$ws.Cells.Item(10,1).Value = "05851583"
$ws.Cells.Item(10,3).Value = "05851583000188"
$ws.Cells.Item(10,8).Value = "31/10/2025 10:38:35"

$ws.Cells.Item(11,1).Value = "08811370"
$ws.Cells.Item(11,3).Value = "08811370000110"
$ws.Cells.Item(11,8).Value = "31/10/2025 10:38:59"

$ws.Cells.Item(12,1).Value = "08811370"
$ws.Cells.Item(12,3).Value = "08811370000200"
$ws.Cells.Item(12,8).Value = "31/10/2025 10:39:12"

$ws.Cells.Item(13,1).Value = "05616807"
$ws.Cells.Item(13,3).Value = "05616807000177"
$ws.Cells.Item(13,8).Value = "31/10/2025 10:39:35"

$ws.Cells.Item(14,1).Value = "37181936"
$ws.Cells.Item(14,3).Value = "37181936000176"
$ws.Cells.Item(14,8).Value = "31/10/2025 10:39:58"

# --- New establishments returned by the latest automation run.

$ws.Cells.Item(15,1).Value = "37181936"
$ws.Cells.Item(15,2).Value = "AGB DISTRIBUIDORA DE ARMAS E MUNICOES S A"
$ws.Cells.Item(15,3).Value = "37181936000257"
$ws.Cells.Item(15,4).Value = "DF"
$ws.Cells.Item(15,5).Value = "BRASILIA"
$ws.Cells.Item(15,6).Value = "2025"
$ws.Cells.Item(15,7).Value = "1,0000"
$ws.Cells.Item(15,8).Value = "31/10/2025 10:40:11"

$ws.Cells.Item(16,1).Value = "07603106"
$ws.Cells.Item(16,2).Value = "AGNUS DEI ARTIGOS RELIGIOSOS LTDA"
$ws.Cells.Item(16,3).Value = "07603106000129"
$ws.Cells.Item(16,4).Value = "GO"
$ws.Cells.Item(16,5).Value = "GOIANIA"
$ws.Cells.Item(16,6).Value = "2025"
$ws.Cells.Item(16,7).Value = "0,5000"
$ws.Cells.Item(16,8).Value = "31/10/2025 10:40:33"
